$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.791.77"
$ws.Range("E2").Value = "  +4.15%  "
$ws.Range("D3").Value = "2.274.09"
$ws.Range("E3").Value = "  +2.17%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'305.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.90%  "
$ws.Range("D6").Value = "'92.92"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.57%  "
$ws.Range("E7").Value = "  +3.89%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +3.68%  "
$ws.Range("D10").Value = "'32.72"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.21%  "
$ws.Range("D11").Value = "'53.82"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.81%  "
$ws.Range("E12").Value = "  +2.56%  "
$ws.Range("E13").Value = "  +1.92%  "
$ws.Range("E14").Value = "  +3.65%  "
$ws.Range("D15").Value = "2.625.87"
$ws.Range("D16").Value = "'14.26"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "2.255.80"
$ws.Range("E17").Value = "  -0.68%  "
$ws.Range("D18").Value = "'0.764"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.32%  "
$ws.Range("D19").Value = "41.731.22"
$ws.Range("E19").Value = "  +4.11%  "
$ws.Range("D20").Value = "'12.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +8.35%  "
$ws.Range("D21").Value = "0.0₃0910"
$ws.Range("E21").Value = "  +2.20%  "
$ws.Range("E22").Value = "  +3.03%  "
$ws.Range("D23").Value = "'67.32"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.41%  "
$ws.Range("D24").Value = "'243.51"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.00%  "
$ws.Range("E25").Value = "  +4.62%  "
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("E27").Value = "  +4.82%  "
$ws.Range("D28").Value = "'24.28"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.49%  "
$ws.Range("D29").Value = "'9.61"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.05%  "
$ws.Range("E30").Value = "  +0.97%  "
$ws.Range("D31").Value = "'34.08"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.84%  "
$ws.Range("D32").Value = "'158.70"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.01%  "
$ws.Range("D33").Value = "'0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.06%  "
$ws.Range("E34").Value = "  +4.11%  "
$ws.Range("D35").Value = "'0.0750"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.55%  "
$ws.Range("E36").Value = "  +0.97%  "
$ws.Range("E37").Value = "  +3.25%  "
$ws.Range("D38").Value = "'16.72"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.46%  "
$ws.Range("E39").Value = "  +5.38%  "
$ws.Range("E40").Value = "  +3.08%  "
$ws.Range("E41").Value = "  +3.67%  "
$ws.Range("E42").Value = "  +5.36%  "
$ws.Range("D43").Value = "2.073.07"
$ws.Range("E43").Value = "  -0.66%  "
$ws.Range("D44").Value = "'19.70"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.60%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "'10.41"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.80%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "'0.0279"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.85%  "
$ws.Range("E47").Value = "  +6.33%  "
$ws.Range("E48").Value = "  +4.67%  "
$ws.Range("E49").Value = "  +3.12%  "
$ws.Range("D50").Value = "'72.94"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.37%  "
$ws.Range("D51").Value = "'1.16"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.73%  "
